$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $cellRef, $val)
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "69.356.87"
$ws.Range("E2").Value = "  +1.59%  "
Set-TextValue $ws "D3" "3.946.25"
$ws.Range("E3").Value = "  +0.53%  "
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue $ws "D5" "492.30"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("E6").Value = "  -0.13%  "
Set-TextValue $ws "D7" "0.625"
$ws.Range("E7").Value = "  -0.29%  "
Set-TextValue $ws "D9" "0.734"
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  +4.53%  "
$ws.Range("E11").Value = "  -1.87%  "
Set-TextValue $ws "D12" "43.18"
$ws.Range("E12").Value = "  +0.89%  "
Set-TextValue $ws "D13" "10.45"
$ws.Range("E13").Value = "  -1.98%  "
Set-TextValue $ws "D14" "4.573.03"
$ws.Range("E14").Value = "  +0.47%  "
Set-TextValue $ws "D15" "3.938.41"
$ws.Range("E15").Value = "  +0.54%  "
Set-TextValue $ws "D16" "14.32"
$ws.Range("E16").Value = "  -3.42%  "
Set-TextValue $ws "D18" "19.90"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("E19").Value = "  +2.22%  "
Set-TextValue $ws "D20" "69.432.64"
$ws.Range("E20").Value = "  +1.57%  "
Set-TextValue $ws "D21" "440.72"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("E22").Value = "  +2.60%  "
Set-TextValue $ws "D23" "14.51"
$ws.Range("E23").Value = "  -1.65%  "
Set-TextValue $ws "D24" "89.10"
$ws.Range("E24").Value = "  +0.65%  "
Set-TextValue $ws "D25" "12.08"
$ws.Range("E25").Value = "  +8.84%  "
$ws.Range("E26").Value = "  +2.43%  "
Set-TextValue $ws "D27" "11.13"
$ws.Range("E27").Value = "  -4.62%  "
Set-TextValue $ws "D28" "37.18"
$ws.Range("E28").Value = "  -4.48%  "
Set-TextValue $ws "D29" "5.62"
$ws.Range("E29").Value = "  -4.29%  "
Set-TextValue $ws "D30" "710.03"
$ws.Range("E30").Value = "  -0.57%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D31" "0.131"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws "D32" "13.36"
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("E33").Value = "  +0.76%  "
Set-TextValue $ws "D34" "0.471"
$ws.Range("E34").Value = "  +25.73%  "
Set-TextValue $ws "D35" "0.0₃0920"
$ws.Range("E35").Value = "  -0.35%  "
Set-TextValue $ws "D36" "61.64"
$ws.Range("E36").Value = "  +4.45%  "
Set-TextValue $ws "D37" "6.05"
$ws.Range("E37").Value = "  +3.91%  "
Set-TextValue $ws "D38" "40.94"
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("E39").Value = "  +0.11%  "
Set-TextValue $ws "D40" "0.999"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +0.04%  "
Set-TextValue $ws "D42" "0.0490"
$ws.Range("E42").Value = "  +2.20%  "
Set-TextValue $ws "D43" "2.95"
$ws.Range("E43").Value = "  +2.44%  "
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("E45").Value = "  +2.21%  "
$ws.Range("E46").Value = "  +0.38%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D47" "3.35"
$ws.Range("E47").Value = "  +7.09%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws "D48" "0.0₆0363"
$ws.Range("E48").Value = "  +5.98%  "
Set-TextValue $ws "D49" "3.06"
$ws.Range("E49").Value = "  +7.61%  "
$ws.Range("E50").Value = "  -0.93%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws "D51" "143.93"
$ws.Range("E51").Value = "  -1.18%  "
